$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "A2" = "三花智控"; "B2" = "闻泰科技"; "C2" = "闻泰科技";
    "A3" = "合肥城建"; "B3" = "白银有色"; "C3" = "道生天合";
    "A4" = "闻泰科技"; "B4" = "三花智控"; "C4" = "三花智控";
    "A5" = "白银有色"; "B5" = "华天科技"; "C5" = "白银有色";
    "A6" = "华天科技"; "B6" = "山子高科"; "C6" = "华天科技";
    "A7" = "山子高科"; "B7" = "东方财富"; "C7" = "山子高科";
    "A8" = "大有能源"; "B8" = "合肥城建"; "C8" = "合肥城建";
    "A9" = "寒武纪-U"; "B9" = "大有能源"; "C9" = "大有能源";
    "A10" = "海峡股份"; "B10" = "寒武纪-U"; "C10" = "楚江新材";
    "A11" = "孚日股份"; "B11" = "贵州茅台"; "C11" = "寒武纪";
    "A12" = "东信和平"; "B12" = "海峡股份"; "C12" = "海峡股份";
    "A13" = "天际股份"; "B13" = "东信和平"; "C13" = "东信和平";
    "A14" = "N道生"; "B14" = "平潭发展"; "C14" = "常山北明";
    "A15" = "北方稀土"; "B15" = "山东墨龙"; "C15" = "远大控股";
    "A16" = "东方财富"; "B16" = "睿能科技"; "C16" = "天际股份";
    "A17" = "海通发展"; "B17" = "国新能源"; "C17" = "安泰科技";
    "A18" = "三孚股份"; "B18" = "安泰集团"; "C18" = "华建集团";
    "A19" = "睿能科技"; "B19" = "三孚股份"; "C19" = "三孚股份";
    "A20" = "建投能源"; "B20" = "孚日股份"; "C20" = "紫金矿业";
    "A21" = "贵州茅台"; "B21" = "紫金矿业"; "C21" = "北方稀土";
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
